$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NSEData")

# Fix header typo in D1: "56 week high" -> "52 week high"
$ws.Range("D1").Value = "52 week high"

# The Face value (E) and 52 week high (D) data columns were swapped for every
# data row. Swap columns D and E back (rows 2-15) using Copy so the original
# cell type (text/shared-string) and formatting are preserved exactly - a
# plain .Value re-assignment would coerce these numeric-looking strings into
# real numbers, which is not what happened here.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 15 }

$dRange = $ws.Range("D2:D" + $lastRow)
$eRange = $ws.Range("E2:E" + $lastRow)
$tmpRange = $ws.Range("Z2:Z" + $lastRow)

$dRange.Copy($tmpRange)
$eRange.Copy($dRange)
$tmpRange.Copy($eRange)
$tmpRange.Clear()
